$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.831.72"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.804.36"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.557"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.86%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.631"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.34"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.53%  "

$ws.Range("E11").Value = "  -2.27%  "

$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.76"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.247.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.809.87"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.949"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.816.21"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.64"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.19"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.72"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.30"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.33"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.46%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.161"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.83"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +9.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.36"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.77%  "

$ws.Range("E31").Value = "  +2.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.16"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "52.32"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.56%  "

$ws.Range("E34").Value = "  +7.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0889"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0446"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.57%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.02"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.65%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.27%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.17"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.37%  "

$ws.Range("E41").Value = "  +1.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.51"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.24"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.89"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.01"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.68%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +8.02%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.42"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.00%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.108.63"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.956"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("E50").Value = "  -1.07%  "

$ws.Range("E51").Value = "  +6.18%  "
